$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.431.26"
$ws.Range("E2").Value = "  +0.35%  "
$ws.Range("D3").Value = "2.624.64"
$ws.Range("E3").Value = "  +0.47%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "'598.36"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.93%  "
$ws.Range("D6").Value = "'153.09"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.89%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  +2.44%  "
$ws.Range("D9").Value = "2.623.90"
$ws.Range("E9").Value = "  +0.48%  "
$ws.Range("E10").Value = "  -1.42%  "
$ws.Range("E11").Value = "  +0.67%  "
$ws.Range("D12").Value = "'5.19"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.56%  "
$ws.Range("D13").Value = "'0.349"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.02%  "
$ws.Range("D14").Value = "'27.70"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.11%  "
$ws.Range("D15").Value = "3.093.80"
$ws.Range("E15").Value = "  +0.17%  "
$ws.Range("D16").Value = "'0.0000182"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.84%  "
$ws.Range("D17").Value = "67.392.47"
$ws.Range("E17").Value = "  +0.61%  "
$ws.Range("D18").Value = "2.618.30"
$ws.Range("E18").Value = "  +0.19%  "
$ws.Range("D19").Value = "'11.14"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.90%  "
$ws.Range("D20").Value = "'363.08"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.00%  "
$ws.Range("D21").Value = "'7.50"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.45%  "
$ws.Range("E22").Value = "  -0.38%  "
$ws.Range("E23").Value = "  +3.87%  "
$ws.Range("E24").Value = "  +0.07%  "
$ws.Range("D25").Value = "'71.08"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +7.36%  "
$ws.Range("E26").Value = "  -1.13%  "
$ws.Range("D27").Value = "2.754.45"
$ws.Range("E27").Value = "  +0.40%  "
$ws.Range("D28").Value = "'586.78"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.79%  "
$ws.Range("B29").Value = "Binance-PegBSC-USD"
$ws.Range("C29").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D29").Value = "'1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.57%  "
$ws.Range("B30").Value = "PEPE"
$ws.Range("C30").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D30").Value = "'0.0000102"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.12%  "
$ws.Range("E31").Value = "  -2.94%  "
$ws.Range("D32").Value = "'7.84"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.63%  "
$ws.Range("D33").Value = "'1.84"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.70%  "
$ws.Range("E34").Value = "  +0.04%  "
$ws.Range("E35").Value = "  -5.51%  "
$ws.Range("E36").Value = "  -1.56%  "
$ws.Range("D37").Value = "'4.90"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.87%  "
$ws.Range("D38").Value = "'157.23"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.43%  "
$ws.Range("D39").Value = "'19.13"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.72%  "
$ws.Range("D40").Value = "'0.369"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.13%  "
$ws.Range("D41").Value = "'5.29"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.29%  "
$ws.Range("E42").Value = "  -0.86%  "
$ws.Range("E43").Value = "  -1.11%  "
$ws.Range("D44").Value = "'41.15"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.05%  "
$ws.Range("E45").Value = "  +0.05%  "
$ws.Range("D46").Value = "'16.35"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.74%  "
$ws.Range("D47").Value = "'156.79"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.60%  "
$ws.Range("D48").Value = "0.0₆0289"
$ws.Range("E48").Value = "  -1.10%  "
$ws.Range("E49").Value = "  -0.15%  "
$ws.Range("B50").Value = "InjectiveProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D50").Value = "'23.13"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +11.39%  "
$ws.Range("B51").Value = "Mantle"
$ws.Range("C51").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D51").Value = "'0.623"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.04%  "
